$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2393.5
$ws.Range("I32").Value = 2274.5
$ws.Range("J32").Value = 2472.8333
$ws.Range("K32").Value = 2274.5
$ws.Range("L32").Value = 2472.8333
$ws.Range("M32").Value = -1948.5
$ws.Range("N32").Value = -3124.8333
$ws.Range("H53").Value = 368.57144
$ws.Range("I53").Value = 313.42856
$ws.Range("K53").Value = 313.42856
$ws.Range("M53").Value = 323.57144
$ws.Range("H64").Value = 9002
$ws.Range("I64").Value = 9002
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 9002
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -8754
$ws.Range("H67").Value = 9002
$ws.Range("I67").Value = 9002
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 9002
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -8144
$ws.Range("N67").ClearContents()
$ws.Range("H76").Value = 6772.846
$ws.Range("I76").Value = 5016.6665
$ws.Range("K76").Value = 5016.6665
$ws.Range("M76").Value = -4701.6665
$ws.Range("H79").Value = 6772.846
$ws.Range("I79").Value = 5016.6665
$ws.Range("K79").Value = 5016.6665
$ws.Range("M79").Value = -3924.6665
$ws.Range("H86").Value = 3663.2058
$ws.Range("J86").Value = 5298.875
$ws.Range("L86").Value = 5298.875
$ws.Range("N86").Value = -7544.875
$ws.Range("H89").Value = 3663.2058
$ws.Range("J89").Value = 5298.875
$ws.Range("L89").Value = 26494.375
$ws.Range("N89").Value = -37726.375
$ws.Range("H116").Value = 182118.1
$ws.Range("I116").Value = 40558.867
$ws.Range("K116").Value = 40558.867
$ws.Range("M116").Value = -37116.867
$ws.Range("H132").Value = 99635.53999999999
$ws.Range("I132").Value = 118465.7
$ws.Range("K132").Value = 355397.1
$ws.Range("M132").Value = -352867.1
$ws.Range("H135").Value = 538.087
$ws.Range("J135").Value = 745.6
$ws.Range("L135").Value = 6710.400000000001
$ws.Range("N135").Value = -11780.4
$ws.Range("H137").Value = 13609.333
$ws.Range("I137").Value = 3898.8
$ws.Range("J137").Value = 25747.5
$ws.Range("K137").Value = 11696.4
$ws.Range("L137").Value = 77242.5
$ws.Range("M137").Value = -9146.400000000001
$ws.Range("N137").Value = -82342.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9673.357
$ws.Range("I2").Value = 4253.263
$ws.Range("K2").Value = 4253.263
$ws.Range("M2").Value = -4140.263
$ws.Range("H45").Value = 3395.45
$ws.Range("I45").Value = 3047.9285
$ws.Range("K45").Value = 3047.9285
$ws.Range("M45").Value = -2670.9285
$ws.Range("H74").Value = 3114.2964
$ws.Range("I74").Value = 1682.2142
$ws.Range("J74").Value = 4656.5386
$ws.Range("K74").Value = 1682.2142
$ws.Range("L74").Value = 4656.5386
$ws.Range("M74").Value = -808.2141999999999
$ws.Range("N74").Value = -6404.5386
$ws.Range("H77").Value = 3114.2964
$ws.Range("I77").Value = 1682.2142
$ws.Range("J77").Value = 4656.5386
$ws.Range("K77").Value = 8411.071
$ws.Range("L77").Value = 23282.693
$ws.Range("M77").Value = -4043.071
$ws.Range("N77").Value = -32018.693
$ws.Range("H110").Value = 1170
$ws.Range("I110").Value = 816.6667
$ws.Range("K110").Value = 816.6667
$ws.Range("M110").Value = 1228.3333
$ws.Range("H116").Value = 9673.357
$ws.Range("I116").Value = 4253.263
$ws.Range("K116").Value = 4253.263
$ws.Range("M116").Value = -1959.263
$ws.Range("H132").Value = 1542761.1
$ws.Range("I132").Value = 2503224.5
$ws.Range("J132").Value = 6020
$ws.Range("K132").Value = 7509673.5
$ws.Range("L132").Value = 18060
$ws.Range("M132").Value = -7507143.5
$ws.Range("N132").Value = -23120

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9673.357
$ws.Range("I3").Value = 4253.263
$ws.Range("K3").Value = 4253.263
$ws.Range("M3").Value = -4139.263
$ws.Range("H22").Value = 305.3846
$ws.Range("I22").Value = 330.66666
$ws.Range("K22").Value = 330.66666
$ws.Range("M22").Value = -157.66666
$ws.Range("H107").Value = 2670.394
$ws.Range("I107").Value = 2077.96
$ws.Range("J107").Value = 4521.75
$ws.Range("K107").Value = 2077.96
$ws.Range("L107").Value = 4521.75
$ws.Range("M107").Value = -157.96
$ws.Range("N107").Value = -8361.75
$ws.Range("H134").Value = 1672586.8
$ws.Range("I134").Value = 1703144.1
$ws.Range("J134").Value = 1529985.4
$ws.Range("K134").Value = 5109432.300000001
$ws.Range("L134").Value = 4589956.199999999
$ws.Range("M134").Value = -5106897.300000001
$ws.Range("N134").Value = -4595026.199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30505.777
$ws.Range("I31").Value = 23496.25
$ws.Range("J31").Value = 32508.5
$ws.Range("K31").Value = 23496.25
$ws.Range("L31").Value = 32508.5
$ws.Range("M31").Value = -23201.25
$ws.Range("N31").Value = -33098.5
$ws.Range("H34").Value = 30505.777
$ws.Range("I34").Value = 23496.25
$ws.Range("J34").Value = 32508.5
$ws.Range("K34").Value = 23496.25
$ws.Range("L34").Value = 32508.5
$ws.Range("M34").Value = -23294.25
$ws.Range("N34").Value = -32912.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 154.6
$ws.Range("I14").Value = 154.6
$ws.Range("K14").Value = 463.8
$ws.Range("M14").Value = -290.8
$ws.Range("H34").Value = 3228.4285
$ws.Range("J34").Value = 5375
$ws.Range("L34").Value = 16125
$ws.Range("N34").Value = -16293
$ws.Range("H39").Value = 1750
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 1750
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("M39").Value = 5250
$ws.Range("N39").Value = -5838
$ws.Range("H55").Value = 2568.2222
$ws.Range("I55").Value = 399.75
$ws.Range("J55").Value = 4303
$ws.Range("K55").Value = 1199.25
$ws.Range("L55").Value = 12909
$ws.Range("M55").Value = -1022.25
$ws.Range("N55").Value = -13263
$ws.Range("H86").Value = 940.8461
$ws.Range("J86").Value = 625
$ws.Range("L86").Value = 1875
$ws.Range("N86").Value = -4247
$ws.Range("H89").Value = 940.8461
$ws.Range("J89").Value = 625
$ws.Range("L89").Value = 5625
$ws.Range("N89").Value = -17481
$ws.Range("H132").Value = 471.66666
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 9502.5
$ws.Range("I18").Value = 9502.5
$ws.Range("K18").Value = 9502.5
$ws.Range("M18").Value = -9209.5
$ws.Range("H70").Value = 3937.647
$ws.Range("I70").Value = 3721.5
$ws.Range("K70").Value = 3721.5
$ws.Range("M70").Value = -3451.5
$ws.Range("H73").Value = 3937.647
$ws.Range("I73").Value = 3721.5
$ws.Range("K73").Value = 3721.5
$ws.Range("M73").Value = -2785.5
$ws.Range("H80").Value = 273189.4
$ws.Range("I80").Value = 369142.16
$ws.Range("K80").Value = 369142.16
$ws.Range("M80").Value = -368144.16
$ws.Range("H83").Value = 273189.4
$ws.Range("I83").Value = 369142.16
$ws.Range("K83").Value = 1845710.8
$ws.Range("M83").Value = -1840718.8
$ws.Range("H102").Value = 2753.6667
$ws.Range("I102").Value = 1826.5769
$ws.Range("J102").Value = 4607.846
$ws.Range("K102").Value = 1826.5769
$ws.Range("L102").Value = 4607.846
$ws.Range("M102").Value = -204.5769
$ws.Range("N102").Value = -7851.846

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2279.9473
$ws.Range("I61").Value = 2015.6428
$ws.Range("K61").Value = 2015.6428
$ws.Range("M61").Value = -1813.6428
$ws.Range("H113").Value = 2279.9473
$ws.Range("I113").Value = 2015.6428
$ws.Range("K113").Value = 2015.6428
$ws.Range("M113").Value = 154.3571999999999
$ws.Range("H122").Value = 4352.9443
$ws.Range("I122").Value = 3926.7097
$ws.Range("K122").Value = 11780.1291
$ws.Range("M122").Value = -9330.1291
$ws.Range("H132").Value = 19860.268
$ws.Range("I132").Value = 26960.5
$ws.Range("J132").Value = 5659.8
$ws.Range("K132").Value = 80881.5
$ws.Range("L132").Value = 16979.4
$ws.Range("M132").Value = -78351.5
$ws.Range("N132").Value = -22039.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 15000
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 15000
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H132").Value = 2489739.5
$ws.Range("I132").Value = 2978174
$ws.Range("K132").Value = 8934522
$ws.Range("M132").Value = -8931992
$ws.Range("H136").Value = 10456.083
$ws.Range("I136").Value = 11568.211
$ws.Range("K136").Value = 34704.633
$ws.Range("M136").Value = -32154.633
